$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.009.50'
$ws.Range("E2").Value = '  +1.82%  '
$ws.Range("E3").Value = '  +1.90%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.59'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.30%  '
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.47'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.95%  '
$ws.Range("E9").Value = '  +1.86%  '
$ws.Range("E10").Value = '  +0.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0872'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("D12").Value = '1.882.83'
$ws.Range("E12").Value = '  +1.94%  '
$ws.Range("D13").Value = '1.647.28'
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.565'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.89%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.66'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("D17").Value = '27.991.34'
$ws.Range("E17").Value = '  +1.84%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.79'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.79%  '
$ws.Range("E19").Value = '  +2.54%  '
$ws.Range("D20").Value = '0.0₃0722'
$ws.Range("E20").Value = '  +0.54%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.68'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.41%  '
$ws.Range("E23").Value = '  +2.72%  '
$ws.Range("E24").Value = '  +4.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.14'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.06%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.69%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E30").Value = '  +1.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0484'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("E32").Value = '  +2.71%  '
$ws.Range("D33").Value = '1.444.71'
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("E34").Value = '  +0.43%  '
$ws.Range("E35").Value = '  +2.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.32'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("E37").Value = '  +3.41%  '
$ws.Range("E38").Value = '  +1.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.560'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.920'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.11%  '
$ws.Range("E42").Value = '  +3.53%  '
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.80'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.40'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.65%  '
$ws.Range("B47").Value = 'MXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.39%  '
$ws.Range("D48").Value = '1.791.57'
$ws.Range("E48").Value = '  +1.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '88.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.68%  '
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("E51").Value = '  +0.28%  '
